$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.568.06'
$ws.Range("E2").Value = '  -3.76%  '
$ws.Range("D3").Value = '2.511.54'
$ws.Range("E3").Value = '  -5.10%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'578.31"
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").Value = "'167.27"
$ws.Range("E6").Value = '  -4.49%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = "'0.517"
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").Value = '2.511.77'
$ws.Range("E9").Value = '  -5.08%  '
$ws.Range("E10").Value = '  -6.73%  '
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("E12").Value = '  -4.11%  '
$ws.Range("D13").Value = "'4.86"
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").Value = '2.964.58'
$ws.Range("E14").Value = '  -5.25%  '
$ws.Range("D15").Value = '69.411.39'
$ws.Range("E15").Value = '  -3.83%  '
$ws.Range("E16").Value = '  -5.59%  '
$ws.Range("D17").Value = "'24.94"
$ws.Range("E17").Value = '  -4.07%  '
$ws.Range("D18").Value = '2.508.09'
$ws.Range("E18").Value = '  -6.15%  '
$ws.Range("D19").Value = "'11.47"
$ws.Range("D20").Value = "'7.77"
$ws.Range("E20").Value = '  -2.95%  '
$ws.Range("D21").Value = "'351.59"
$ws.Range("E21").Value = '  -4.93%  '
$ws.Range("D22").Value = "'3.96"
$ws.Range("E22").Value = '  -4.76%  '
$ws.Range("D23").Value = "'2.00"
$ws.Range("E23").Value = '  -3.14%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = "'69.19"
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("E26").Value = '  -5.58%  '
$ws.Range("D27").Value = "'9.06"
$ws.Range("E27").Value = '  -6.47%  '
$ws.Range("D28").Value = '2.641.44'
$ws.Range("E28").Value = '  -5.04%  '
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").Value = '0.0₃0906'
$ws.Range("E30").Value = '  -5.14%  '
$ws.Range("D31").Value = "'7.89"
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("D32").Value = "'478.72"
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("E34").Value = '  -2.89%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = "'153.91"
$ws.Range("E37").Value = '  -5.57%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = "'18.60"
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").Value = "'4.78"
$ws.Range("E41").Value = '  -2.78%  '
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("E43").Value = '  -6.67%  '
$ws.Range("E44").Value = '  -13.58%  '
$ws.Range("E45").Value = '  -8.49%  '
$ws.Range("D46").Value = "'38.21"
$ws.Range("E46").Value = '  -2.46%  '
$ws.Range("D47").Value = "'144.39"
$ws.Range("E47").Value = '  -6.20%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = "'0.532"
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = "'3.55"
$ws.Range("E49").Value = '  -3.39%  '
$ws.Range("E50").Value = '  -4.95%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0730"
$ws.Range("E51").Value = '  -2.39%  '

# Clear formatting on cells that needed quote-prefix to stay text,
# so no stray number-format style is introduced.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D51").ClearFormats()
